# Add one previously-missed OpinionWay poll (rows 44-46) and one previously-
# missed Ifop poll (rows 47-52), and correct the sample-size (n) column for
# several existing rows whose "included" n had mistakenly been left equal to
# the "partially included" n (plus two week-of-fieldwork corrections on
# E17/E18) - per commit "add one ifop poll and one opinionway poll
# (previously missed)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix sample sizes (column I) on existing rows 14-43, and the two
#    fieldwork-day corrections on E17/E18.
# ---------------------------------------------------------------------
$ws.Range("I14").Value = 891
$ws.Range("I15").Value = 880
$ws.Range("I16").Value = 859

$ws.Range("E17").Value = 19
$ws.Range("I17").Value = 869

$ws.Range("E18").Value = 19
$ws.Range("I18").Value = 838

$ws.Range("I19").Value = 872
$ws.Range("I20").Value = 861
$ws.Range("I21").Value = 850
$ws.Range("I22").Value = 893

$ws.Range("I23").Value = 904
$ws.Range("I24").Value = 915
$ws.Range("I25").Value = 872
$ws.Range("I26").Value = 904
$ws.Range("I27").Value = 904

$ws.Range("I29").Value = 624
$ws.Range("I30").Value = 642

$ws.Range("I31").Value = 825
$ws.Range("I32").Value = 806
$ws.Range("I33").Value = 815
$ws.Range("I34").Value = 815

$ws.Range("I35").Value = 956
$ws.Range("I36").Value = 931
$ws.Range("I37").Value = 920

$ws.Range("I38").Value = 807
$ws.Range("I39").Value = 788
$ws.Range("I40").Value = 779

$ws.Range("I41").Value = 892
$ws.Range("I42").Value = 903
$ws.Range("I43").Value = 892

# ---------------------------------------------------------------------
# 2) Append the missed OpinionWay poll (week 13, 2021-09-23), 3 rows.
# ---------------------------------------------------------------------
$ws.Range("A44").Value = 13
$ws.Range("B44").Value = 2021
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 9
$ws.Range("E44").Value = 23
$ws.Range("F44").Value = "opinionway"
$ws.Range("G44").Value = "online"
$ws.Range("H44").Value = "partially"
$ws.Range("I44").Value = 709
$ws.Range("J44").Value = 1
$ws.Range("K44").Value = 1
$ws.Range("L44").Value = 7.5
$ws.Range("M44").Value = 2
$ws.Range("N44").Value = 2
$ws.Range("O44").Value = 9
$ws.Range("P44").Value = 6
$ws.Range("Q44").Value = 26.5
$ws.Range("T44").Value = 12.5
$ws.Range("V44").Value = 2
$ws.Range("W44").Value = 20
$ws.Range("X44").Value = 10.5

$ws.Range("A45").Value = 13
$ws.Range("B45").Value = 2021
$ws.Range("C45").Value = 4
$ws.Range("D45").Value = 9
$ws.Range("E45").Value = 23
$ws.Range("F45").Value = "opinionway"
$ws.Range("G45").Value = "online"
$ws.Range("H45").Value = "partially"
$ws.Range("I45").Value = 692
$ws.Range("J45").Value = 1.5
$ws.Range("K45").Value = 1
$ws.Range("L45").Value = 7
$ws.Range("M45").Value = 2
$ws.Range("N45").Value = 2
$ws.Range("O45").Value = 9.5
$ws.Range("P45").Value = 7
$ws.Range("Q45").Value = 28
$ws.Range("R45").Value = 12
$ws.Range("V45").Value = 4.5
$ws.Range("W45").Value = 25.5

$ws.Range("A46").Value = 13
$ws.Range("B46").Value = 2021
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = 9
$ws.Range("E46").Value = 23
$ws.Range("F46").Value = "opinionway"
$ws.Range("G46").Value = "online"
$ws.Range("H46").Value = "partially"
$ws.Range("I46").Value = 689
$ws.Range("J46").Value = 1.5
$ws.Range("K46").Value = 1
$ws.Range("L46").Value = 8
$ws.Range("M46").Value = 2.5
$ws.Range("N46").Value = 2.5
$ws.Range("O46").Value = 9
$ws.Range("P46").Value = 7
$ws.Range("Q46").Value = 29.5
$ws.Range("S46").Value = 7.5
$ws.Range("V46").Value = 5.5
$ws.Range("W46").Value = 26.5

# ---------------------------------------------------------------------
# 3) Append the missed Ifop poll (week 14, 2021-09-29), 6 rows.
# ---------------------------------------------------------------------
$ws.Range("A47").Value = 14
$ws.Range("B47").Value = 2021
$ws.Range("C47").Value = 5
$ws.Range("D47").Value = 9
$ws.Range("E47").Value = 29
$ws.Range("F47").Value = "ifop"
$ws.Range("G47").Value = "online"
$ws.Range("H47").Value = "included"
$ws.Range("I47").Value = 959
$ws.Range("J47").Value = 1
$ws.Range("K47").Value = 0.5
$ws.Range("L47").Value = 7.5
$ws.Range("M47").Value = 2
$ws.Range("N47").Value = 4
$ws.Range("O47").Value = 8.5
$ws.Range("P47").Value = 6
$ws.Range("Q47").Value = 26
$ws.Range("T47").Value = 18
$ws.Range("U47").Value = 1
$ws.Range("V47").Value = 2.5
$ws.Range("W47").Value = 23

$ws.Range("A48").Value = 14
$ws.Range("B48").Value = 2021
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 9
$ws.Range("E48").Value = 29
$ws.Range("F48").Value = "ifop"
$ws.Range("G48").Value = "online"
$ws.Range("H48").Value = "included"
$ws.Range("I48").Value = 939
$ws.Range("J48").Value = 1
$ws.Range("K48").Value = 0.5
$ws.Range("L48").Value = 8
$ws.Range("M48").Value = 2.5
$ws.Range("N48").Value = 4
$ws.Range("O48").Value = 8.5
$ws.Range("P48").Value = 6
$ws.Range("Q48").Value = 27
$ws.Range("R48").Value = 13
$ws.Range("U48").Value = 1
$ws.Range("V48").Value = 3.5
$ws.Range("W48").Value = 25

$ws.Range("A49").Value = 14
$ws.Range("B49").Value = 2021
$ws.Range("C49").Value = 5
$ws.Range("D49").Value = 9
$ws.Range("E49").Value = 29
$ws.Range("F49").Value = "ifop"
$ws.Range("G49").Value = "online"
$ws.Range("H49").Value = "included"
$ws.Range("I49").Value = 934
$ws.Range("J49").Value = 1
$ws.Range("K49").Value = 0.5
$ws.Range("L49").Value = 8
$ws.Range("M49").Value = 2
$ws.Range("N49").Value = 3.5
$ws.Range("O49").Value = 9
$ws.Range("P49").Value = 6
$ws.Range("Q49").Value = 27
$ws.Range("S49").Value = 12
$ws.Range("U49").Value = 2
$ws.Range("V49").Value = 4
$ws.Range("W49").Value = 25

$ws.Range("A50").Value = 14
$ws.Range("B50").Value = 2021
$ws.Range("C50").Value = 5
$ws.Range("D50").Value = 9
$ws.Range("E50").Value = 29
$ws.Range("F50").Value = "ifop"
$ws.Range("G50").Value = "online"
$ws.Range("H50").Value = "included"
$ws.Range("I50").Value = 990
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = 0.5
$ws.Range("L50").Value = 7
$ws.Range("M50").Value = 2
$ws.Range("N50").Value = 3.5
$ws.Range("O50").Value = 8
$ws.Range("P50").Value = 5.5
$ws.Range("Q50").Value = 24
$ws.Range("T50").Value = 15
$ws.Range("U50").Value = 1
$ws.Range("V50").Value = 2.5
$ws.Range("W50").Value = 18
$ws.Range("X50").Value = 12

$ws.Range("A51").Value = 14
$ws.Range("B51").Value = 2021
$ws.Range("C51").Value = 5
$ws.Range("D51").Value = 9
$ws.Range("E51").Value = 29
$ws.Range("F51").Value = "ifop"
$ws.Range("G51").Value = "online"
$ws.Range("H51").Value = "included"
$ws.Range("I51").Value = 976
$ws.Range("J51").Value = 1
$ws.Range("K51").Value = 0.5
$ws.Range("L51").Value = 7
$ws.Range("M51").Value = 2
$ws.Range("N51").Value = 3
$ws.Range("O51").Value = 8
$ws.Range("P51").Value = 6
$ws.Range("Q51").Value = 26
$ws.Range("R51").Value = 11
$ws.Range("U51").Value = 1
$ws.Range("V51").Value = 2.5
$ws.Range("W51").Value = 19
$ws.Range("X51").Value = 13

$ws.Range("A52").Value = 14
$ws.Range("B52").Value = 2021
$ws.Range("C52").Value = 5
$ws.Range("D52").Value = 9
$ws.Range("E52").Value = 29
$ws.Range("F52").Value = "ifop"
$ws.Range("G52").Value = "online"
$ws.Range("H52").Value = "included"
$ws.Range("I52").Value = 967
$ws.Range("J52").Value = 1
$ws.Range("K52").Value = 0.5
$ws.Range("L52").Value = 7
$ws.Range("M52").Value = 1.5
$ws.Range("N52").Value = 3
$ws.Range("O52").Value = 8.5
$ws.Range("P52").Value = 6
$ws.Range("Q52").Value = 26
$ws.Range("S52").Value = 9.5
$ws.Range("U52").Value = 1
$ws.Range("V52").Value = 2.5
$ws.Range("W52").Value = 19.5
$ws.Range("X52").Value = 14

# ---------------------------------------------------------------------
# 4) Leave the final selection on the last entered cell, matching the
#    author's post-edit cursor position.
# ---------------------------------------------------------------------
$ws.Range("E52").Select()
